$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 into the new
# header cells I1 and J1, then set their text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-35.
$values = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(1, 1),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(5, 5),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(9, 9),
    @(3, 4),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
